$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:G2").NumberFormat = "@"
$ws.Range("D2").Value = '319.52'
$ws.Range("E2").Value = '1.43%'
$ws.Range("F2").Value = '2-2-2023'
$ws.Range("G2").Value = '0'

$ws.Range("D3:G3").NumberFormat = "@"
$ws.Range("D3").Value = '40.07'
$ws.Range("E3").Value = '5.85%'
$ws.Range("F3").Value = '2-2-2023'
$ws.Range("G3").Value = '0'

$ws.Range("D4:G4").NumberFormat = "@"
$ws.Range("D4").Value = '5.148'
$ws.Range("E4").Value = '0.27%'
$ws.Range("F4").Value = '2-2-2023'
$ws.Range("G4").Value = '0'

$ws.Range("D5:G5").NumberFormat = "@"
$ws.Range("D5").Value = '0.08148'
$ws.Range("E5").Value = '3.34%'
$ws.Range("F5").Value = '2-2-2023'
$ws.Range("G5").Value = '0'

$ws.Range("D6:G6").NumberFormat = "@"
$ws.Range("B6").Value = 'KuCoinToken'
$ws.Range("C6").Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
$ws.Range("D6").Value = '8.573'
$ws.Range("E6").Value = '3.67%'
$ws.Range("F6").Value = '2-2-2023'
$ws.Range("G6").Value = '0'

$ws.Range("D7:G7").NumberFormat = "@"
$ws.Range("B7").Value = 'FTXToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D7").Value = '1.917'
$ws.Range("E7").Value = '0.67%'
$ws.Range("F7").Value = '2-2-2023'
$ws.Range("G7").Value = '0'

$ws.Range("D8:G8").NumberFormat = "@"
$ws.Range("B8").Value = 'BTSEToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D8").Value = '2.954'
$ws.Range("E8").Value = '-1.88%'
$ws.Range("F8").Value = '2-2-2023'
$ws.Range("G8").Value = '0'

$ws.Range("D9:G9").NumberFormat = "@"
$ws.Range("B9").Value = 'MXToken'
$ws.Range("C9").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D9").Value = '0.9449'
$ws.Range("E9").Value = '2.58%'
$ws.Range("F9").Value = '2-2-2023'
$ws.Range("G9").Value = '0'

$ws.Range("D10:G10").NumberFormat = "@"
$ws.Range("B10").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C10").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D10").Value = '0.1297'
$ws.Range("E10").Value = '14.68%'
$ws.Range("F10").Value = '2-2-2023'
$ws.Range("G10").Value = '0'

$ws.Range("D11:G11").NumberFormat = "@"
$ws.Range("B11").Value = 'WazirX'
$ws.Range("C11").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D11").Value = '0.1958'
$ws.Range("E11").Value = '1.32%'
$ws.Range("F11").Value = '2-2-2023'
$ws.Range("G11").Value = '0'

$ws.Range("D12:G12").NumberFormat = "@"
$ws.Range("B12").Value = 'MandalaExchangeToken'
$ws.Range("C12").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D12").Value = '0.09065'
$ws.Range("E12").Value = '-0.31%'
$ws.Range("F12").Value = '2-2-2023'
$ws.Range("G12").Value = '0'

$ws.Range("D13:G13").NumberFormat = "@"
$ws.Range("B13").Value = 'BitrueCoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D13").Value = '0.03452'
$ws.Range("E13").Value = '3.01%'
$ws.Range("F13").Value = '2-2-2023'
$ws.Range("G13").Value = '0'

$ws.Range("D14:G14").NumberFormat = "@"
$ws.Range("B14").Value = 'BitMartToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D14").Value = '0.09549'
$ws.Range("E14").Value = '-0.59%'
$ws.Range("F14").Value = '2-2-2023'
$ws.Range("G14").Value = '0'

$ws.Range("D15:G15").NumberFormat = "@"
$ws.Range("B15").Value = 'BitForexToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D15").Value = '0.001406'
$ws.Range("E15").Value = '1.80%'
$ws.Range("F15").Value = '2-2-2023'
$ws.Range("G15").Value = '0'

$ws.Range("D16:G16").NumberFormat = "@"
$ws.Range("B16").Value = 'TigerCash'
$ws.Range("C16").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D16").Value = '0.005911'
$ws.Range("E16").Value = '-1.44%'
$ws.Range("F16").Value = '2-2-2023'
$ws.Range("G16").Value = '0'

$ws.Range("D17:G17").NumberFormat = "@"
$ws.Range("B17").Value = 'LEO'
$ws.Range("C17").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D17").Value = '3.360'
$ws.Range("E17").Value = '-5.37%'
$ws.Range("F17").Value = '2-2-2023'
$ws.Range("G17").Value = '0'

$ws.Range("D18:G18").NumberFormat = "@"
$ws.Range("B18").Value = 'GateToken'
$ws.Range("C18").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D18").Value = '4.468'
$ws.Range("E18").Value = '1.04%'
$ws.Range("F18").Value = '2-2-2023'
$ws.Range("G18").Value = '0'

$ws.Range("D19:G19").NumberFormat = "@"
$ws.Range("D19").Value = '0.3534'
$ws.Range("E19").Value = '2.67%'
$ws.Range("F19").Value = '2-2-2023'
$ws.Range("G19").Value = '0'

$ws.Range("D20:G20").NumberFormat = "@"
$ws.Range("D20").Value = '6.579'
$ws.Range("E20").Value = '24.60%'
$ws.Range("F20").Value = '2-2-2023'
$ws.Range("G20").Value = '0'

$ws.Range("D21:G21").NumberFormat = "@"
$ws.Range("D21").Value = '0.1324'
$ws.Range("E21").Value = '2.74%'
$ws.Range("F21").Value = '2-2-2023'
$ws.Range("G21").Value = '0'

$ws.Range("D22:G22").NumberFormat = "@"
$ws.Range("D22").Value = '0.2305'
$ws.Range("E22").Value = '-11.17%'
$ws.Range("F22").Value = '2-2-2023'
$ws.Range("G22").Value = '0'

$ws.Range("D23:G23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04388'
$ws.Range("E23").Value = '0.44%'
$ws.Range("F23").Value = '2-2-2023'
$ws.Range("G23").Value = '0'

$ws.Range("D24:G24").NumberFormat = "@"
$ws.Range("D24").Value = '0.001222'
$ws.Range("E24").Value = '-1.98%'
$ws.Range("F24").Value = '2-2-2023'
$ws.Range("G24").Value = '0'

$ws.Range("D25:G25").NumberFormat = "@"
$ws.Range("D25").Value = '0.004330'
$ws.Range("E25").Value = '-7.31%'
$ws.Range("F25").Value = '2-2-2023'
$ws.Range("G25").Value = '0'

$ws.Range("D26:G26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0001129'
$ws.Range("E26").Value = '-17.04%'
$ws.Range("F26").Value = '2-2-2023'
$ws.Range("G26").Value = '0'

$ws.Range("D27:G27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0003978'
$ws.Range("E27").Value = '-0.50%'
$ws.Range("F27").Value = '2-2-2023'
$ws.Range("G27").Value = '0'

$ws.Range("F28:G28").NumberFormat = "@"
$ws.Range("F28").Value = '2-2-2023'
$ws.Range("G28").Value = '0'

$ws.Range("F29:G29").NumberFormat = "@"
$ws.Range("F29").Value = '2-2-2023'
$ws.Range("G29").Value = '0'

$ws.Range("F30:G30").NumberFormat = "@"
$ws.Range("F30").Value = '2-2-2023'
$ws.Range("G30").Value = '0'

$ws.Range("F31:G31").NumberFormat = "@"
$ws.Range("F31").Value = '2-2-2023'
$ws.Range("G31").Value = '0'

$ws.Range("F32:G32").NumberFormat = "@"
$ws.Range("F32").Value = '2-2-2023'
$ws.Range("G32").Value = '0'

$ws.Range("F33:G33").NumberFormat = "@"
$ws.Range("F33").Value = '2-2-2023'
$ws.Range("G33").Value = '0'

$ws.Range("F34:G34").NumberFormat = "@"
$ws.Range("F34").Value = '2-2-2023'
$ws.Range("G34").Value = '0'

$ws.Range("F35:G35").NumberFormat = "@"
$ws.Range("F35").Value = '2-2-2023'
$ws.Range("G35").Value = '0'

$ws.Range("F36:G36").NumberFormat = "@"
$ws.Range("F36").Value = '2-2-2023'
$ws.Range("G36").Value = '0'

$ws.Range("F37:G37").NumberFormat = "@"
$ws.Range("F37").Value = '2-2-2023'
$ws.Range("G37").Value = '0'

$ws.Range("F38:G38").NumberFormat = "@"
$ws.Range("F38").Value = '2-2-2023'
$ws.Range("G38").Value = '0'

$ws.Range("D39:G39").NumberFormat = "@"
$ws.Range("D39").Value = '0.02372'
$ws.Range("E39").Value = '4.31%'
$ws.Range("F39").Value = '2-2-2023'
$ws.Range("G39").Value = '0'

$ws.Range("D40:G40").NumberFormat = "@"
$ws.Range("D40").Value = '0.05207'
$ws.Range("E40").Value = '1.76%'
$ws.Range("F40").Value = '2-2-2023'
$ws.Range("G40").Value = '0'

$ws.Range("E41:G41").NumberFormat = "@"
$ws.Range("E41").Value = '3.00%'
$ws.Range("F41").Value = '2-2-2023'
$ws.Range("G41").Value = '0'

$ws.Range("D42:G42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1424'
$ws.Range("E42").Value = '5.03%'
$ws.Range("F42").Value = '2-2-2023'
$ws.Range("G42").Value = '0'

$ws.Range("D43:G43").NumberFormat = "@"
$ws.Range("D43").Value = '0.008621'
$ws.Range("E43").Value = '-4.62%'
$ws.Range("F43").Value = '2-2-2023'
$ws.Range("G43").Value = '0'

$ws.Range("D44:G44").NumberFormat = "@"
$ws.Range("D44").Value = '0.002028'
$ws.Range("E44").Value = '3.93%'
$ws.Range("F44").Value = '2-2-2023'
$ws.Range("G44").Value = '0'

$ws.Range("D45:G45").NumberFormat = "@"
$ws.Range("D45").Value = '0.009191'
$ws.Range("E45").Value = '6.38%'
$ws.Range("F45").Value = '2-2-2023'
$ws.Range("G45").Value = '0'

$ws.Range("D46:G46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00006431'
$ws.Range("E46").Value = '-3.25%'
$ws.Range("F46").Value = '2-2-2023'
$ws.Range("G46").Value = '0'

$ws.Range("D47:G47").NumberFormat = "@"
$ws.Range("D47").Value = '0.00000000748'
$ws.Range("E47").Value = '-0.36%'
$ws.Range("F47").Value = '2-2-2023'
$ws.Range("G47").Value = '0'

$ws.Range("D48:G48").NumberFormat = "@"
$ws.Range("D48").Value = '0.002834'
$ws.Range("E48").Value = '-11.85%'
$ws.Range("F48").Value = '2-2-2023'
$ws.Range("G48").Value = '0'

$ws.Range("D49:G49").NumberFormat = "@"
$ws.Range("D49").Value = '0.001685'
$ws.Range("E49").Value = '68.09%'
$ws.Range("F49").Value = '2-2-2023'
$ws.Range("G49").Value = '0'

$ws.Range("D50:G50").NumberFormat = "@"
$ws.Range("D50").Value = '0.00002094'
$ws.Range("E50").Value = '-0.36%'
$ws.Range("F50").Value = '2-2-2023'
$ws.Range("G50").Value = '0'

$ws.Range("D51:G51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0001994'
$ws.Range("E51").Value = '-0.36%'
$ws.Range("F51").Value = '2-2-2023'
$ws.Range("G51").Value = '0'

